$wb = $excel.ActiveWorkbook

# Sheet "展览" - column F ("想去人数") updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 56
$wsExpo.Range("F3").Value = 785
$wsExpo.Range("F4").Value = 39
$wsExpo.Range("F6").Value = 71
$wsExpo.Range("F7").Value = 272
$wsExpo.Range("F8").Value = 3911
$wsExpo.Range("F9").Value = 90
$wsExpo.Range("F10").Value = 4600
$wsExpo.Range("F11").Value = 505
$wsExpo.Range("F12").Value = 1160
$wsExpo.Range("F13").Value = 72

# Sheet "全部类型" - column F ("想去人数") updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 56
$wsAll.Range("F3").Value = 785
$wsAll.Range("F4").Value = 39
$wsAll.Range("F6").Value = 71
$wsAll.Range("F8").Value = 272
$wsAll.Range("F9").Value = 3911
$wsAll.Range("F10").Value = 90
$wsAll.Range("F11").Value = 4600
$wsAll.Range("F12").Value = 505
$wsAll.Range("F13").Value = 1160
$wsAll.Range("F14").Value = 72
